{"js": "// Insert the four new \"attribute\" lines before the existing (last) paragraph,\n// which holds the _GoBack bookmark and must remain the final paragraph.\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nconst firstParagraph = paragraphs.items[0];\n\n// Insert three new paragraphs before the original (bookmark-holding) one.\nfirstParagraph.insertParagraph(\"String nombre;\", Word.InsertLocation.before);\nfirstParagraph.insertParagraph(\"int edad;\", Word.InsertLocation.before);\nfirstParagraph.insertParagraph(\"Date fecha;\", Word.InsertLocation.before);\n\n// The original paragraph (still holding the _GoBack bookmark) becomes the\n// \"long cedula;\" line -- insert the text at its start rather than adding a\n// brand new paragraph, so the bookmark stays in the same <w:p>.\nfirstParagraph.insertText(\"long cedula;\", Word.InsertLocation.start);\n\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# The document starts with exactly one (empty) paragraph that holds the\n# _GoBack bookmark. Add three new empty paragraphs in front of it, fill in\n# the first three with their text, then prepend the fourth line's text\n# directly into the original (bookmark-holding) paragraph so the bookmark\n# stays attached to the \"long cedula;\" line, matching the authored edit.\n$firstPara = $d.Paragraphs.First\n$r = $firstPara.Range\n$r.InsertParagraphBefore()\n$r.InsertParagraphBefore()\n$r.InsertParagraphBefore()\n\n$d.Paragraphs.Item(1).Range.Text = \"String nombre;\"\n$d.Paragraphs.Item(2).Range.Text = \"int edad;\"\n$d.Paragraphs.Item(3).Range.Text = \"Date fecha;\"\n\n$lastPara = $d.Paragraphs.Item($d.Paragraphs.Count)\n$lastPara.Range.InsertBefore(\"long cedula;\")\n\n$d.Save()\n"}
